# Insert two new weekly price-report rows (a new "week"/date group) at
# row 319 of Sheet1, pushing the existing rows 319:423 down to 321:425.
# This mirrors the commit "Fruta / hortaliza, semanal" (weekly update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 319 (data below shifts down by 2).
$ws.Rows.Item(319).Insert()
$ws.Rows.Item(319).Insert()

# --- New row 319 -----------------------------------------------------
$ws.Cells.Item(319, 1).Value2 = 5
$ws.Cells.Item(319, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(319, 3).Value2 = "Maule"
$ws.Cells.Item(319, 4).Value2 = 44524
$ws.Cells.Item(319, 5).Value2 = 7
$ws.Cells.Item(319, 6).Value2 = 100112020
$ws.Cells.Item(319, 7).Value2 = "Tomate"
$ws.Cells.Item(319, 8).Value2 = "Larga vida"
$ws.Cells.Item(319, 9).Value2 = "Primera"
$ws.Cells.Item(319, 10).Value2 = 2000
$ws.Cells.Item(319, 11).Value2 = 13000
$ws.Cells.Item(319, 12).Value2 = 13000
$ws.Cells.Item(319, 13).Value2 = 13000
$ws.Cells.Item(319, 14).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(319, 15).Value2 = "Región del Maule"
$ws.Cells.Item(319, 16).Value2 = 722
$ws.Cells.Item(319, 17).Value2 = 18
$ws.Cells.Item(319, 18).Value2 = "Hortaliza"

# --- New row 320 -----------------------------------------------------
$ws.Cells.Item(320, 1).Value2 = 5
$ws.Cells.Item(320, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(320, 3).Value2 = "Maule"
$ws.Cells.Item(320, 4).Value2 = 44524
$ws.Cells.Item(320, 5).Value2 = 7
$ws.Cells.Item(320, 6).Value2 = 100112020
$ws.Cells.Item(320, 7).Value2 = "Tomate"
$ws.Cells.Item(320, 8).Value2 = "Larga vida"
$ws.Cells.Item(320, 9).Value2 = "Primera"
$ws.Cells.Item(320, 10).Value2 = 3500
$ws.Cells.Item(320, 11).Value2 = 8000
$ws.Cells.Item(320, 12).Value2 = 8000
$ws.Cells.Item(320, 13).Value2 = 8000
$ws.Cells.Item(320, 14).Value2 = "`$/caja 15 kilos"
$ws.Cells.Item(320, 15).Value2 = "Región del Maule"
$ws.Cells.Item(320, 16).Value2 = 533
$ws.Cells.Item(320, 17).Value2 = 15
$ws.Cells.Item(320, 18).Value2 = "Hortaliza"
